$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6425934433937073
$ws.Range("B1").Value = 0.7231345176696777
$ws.Range("C1").Value = 0.8819990158081055
$ws.Range("D1").Value = 1.662474393844604
$ws.Range("E1").Value = 5.366727352142334
